$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Find-ParaIndexByText($searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($searchText)) {
            return $i
        }
    }
    return -1
}

# Bold only the leading $termLen characters of paragraph $paraIdx (a
# "Term: description" style paragraph), leaving the remainder unbolded.
function Bold-Term($paraIdx, $termLen, $totalLen) {
    $p = $d.Paragraphs.Item($paraIdx).Range
    $p.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
    $start = $p.Start
    $p.Bold = $true
    if ($totalLen -gt $termLen) {
        $tail = $d.Range($start + $termLen, $start + $totalLen)
        $tail.Bold = $false
    }
}

# Insert a new "Term: description" glossary paragraph immediately before
# the paragraph currently at $beforeParaIndex. Returns the index of the
# freshly created paragraph.
function Insert-DefParagraphBefore($beforeParaIndex, $term, $desc) {
    $target = $d.Paragraphs.Item($beforeParaIndex).Range
    $target.Collapse(1) | Out-Null   # wdCollapseStart
    $full = $term + $desc
    $target.InsertBefore($full + "`r")
    Bold-Term $beforeParaIndex $term.Length $full.Length
    return $beforeParaIndex
}

# Insert a new "Term: description" glossary paragraph immediately after
# the paragraph currently at $afterParaIndex. Returns the index of the
# freshly created paragraph.
function Insert-DefParagraphAfter($afterParaIndex, $term, $desc) {
    $target = $d.Paragraphs.Item($afterParaIndex).Range
    $target.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
    $target.Collapse(0) | Out-Null      # wdCollapseEnd
    $full = $term + $desc
    $target.InsertAfter("`r" + $full)
    $newIdx = $afterParaIndex + 1
    Bold-Term $newIdx $term.Length $full.Length
    return $newIdx
}

# ---------------------------------------------------------------------------
# 1. Algorithm  -- new paragraph, right before "Application Program
#    Interface (API)"
# ---------------------------------------------------------------------------

$apiIdx = Find-ParaIndexByText "Application Program Interface (API)"
Insert-DefParagraphBefore $apiIdx "Algorithm" ": a set of instructions for solving a problem" | Out-Null

# ---------------------------------------------------------------------------
# 2. Concatenation -- new paragraph, right after "Application Program
#    Interface (API)"; also carries the relocated "_GoBack" bookmark
#    (previously sitting after the IDE paragraph). Adding a bookmark
#    named "_GoBack" automatically removes any pre-existing one elsewhere
#    in the document, which reproduces the diff's move.
# ---------------------------------------------------------------------------

$apiIdx = Find-ParaIndexByText "Application Program Interface (API)"
$concatTerm = "Concatenation"
$concatDesc = ": the process of combining things together like a pair of strings or several lists"
$concatIdx = Insert-DefParagraphAfter $apiIdx $concatTerm $concatDesc

$concatPara = $d.Paragraphs.Item($concatIdx).Range
$concatPara.MoveEnd(1, -1) | Out-Null
$concatStart = $concatPara.Start
$beforeBookmark = $concatTerm + ": the process of combining things together like a pair of strings or several"
$bmPos = $concatStart + $beforeBookmark.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. Compiler -- new paragraph, right before "Documentation"
# ---------------------------------------------------------------------------

$docIdx = Find-ParaIndexByText "Documentation"
Insert-DefParagraphBefore $docIdx "Compiler" ": a software system which converts source code to lower-level code" | Out-Null

# ---------------------------------------------------------------------------
# 4. Variable -- drop the trailing space in its description
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    ": a value that is subject to change ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ": a value that is subject to change", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Version Control -- new paragraph, right after "Variable"
# ---------------------------------------------------------------------------

$varIdx = Find-ParaIndexByText "Variable"
Insert-DefParagraphAfter $varIdx "Version Control" ": a system that supports organization of many versions of software" | Out-Null
